$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only these columns actually differ between row 2 and row 3 in the source
# data; swap them cell-by-cell. (Columns such as Y/AA hold identical
# "2026-01-24" text in both rows, so leaving them untouched both matches
# the diff and avoids Excel's automatic text->date coercion on write.)
$cols = "A","B","D","E","F","G","H","K","Q","R","AH","AJ","AK","AO"

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
